# Fix Training Data Issue (#48)
# The "Date" column (BF) held a malformed literal "6-7-2007-08" for every
# data row. It should read the correct ISO date text "2008-06-07".
#
# NOTE: assigning a date-shaped string straight to .Value lets the host
# auto-detect it as a real date serial, which is not what we want here -
# the column must keep holding plain text. Forcing the cell to Text format
# before the write keeps it literal; flipping the style back to "Normal"
# afterwards removes the now-unneeded explicit formatting so the cell's
# style matches its original (default) styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 31
$col = 58  # column BF

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $cell.NumberFormat = "@"
    $cell.Value = "2008-06-07"
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}
